$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (Location | Site | Exposure period | Notes) for rows 2-17.
# Two old rows (Bermagui NSW, Eden NSW) were removed and a new
# "Lakes Boat Shed Cafe" row was inserted among the Lakes Entrance entries,
# shrinking the table from 17 data rows (A1:D18) to 16 data rows (A1:D17).
$data = @(
    @('Brighton', 'Sons of Mary Restaurant  14 Spink St, Brighton VIC 3186', '24/12/2020 10:00am-11:05am', 'Case ate at restaurant'),
    @('Doveton', 'Holy Family Parish Doveton Catholic  100 Power Road, Doveton VIC 3177', '26/12/20 6:30pm', 'Case attended English service'),
    @('Lakes Entrance', 'Albert and Co. Cafe - Bellevue Hotel  201 Esplanade, Lakes Entrance VIC 3909', '29/12/2020 09:30am-10:45am', 'Case ate in cafe'),
    @('Lakes Entrance', 'Esplanade Resort Lakes Entrance - Hotel bar  1 Esplanade, Lakes Entrance VIC 3909', '29/12/2020 05:30pm-05:50pm', 'Case in hotel bar'),
    @('Lakes Entrance', 'Esplanade Resort Lakes Entrance - Pool area  1 Esplanade, Lakes Entrance VIC 3909', '29/12/2020 03:30pm-04:30pm', 'Case visited pool area'),
    @('Lakes Entrance', 'Lakes Boat Shed Cafe  54 Marine Parade, Lakes Entrance VIC 3909', '30/12/2020 9:30am-10:30am', 'Case attended cafe'),
    @('Lakes Entrance', 'The Esplanade Resort and Spa  1 Esplanade, Lakes Entrance VIC 3909', '29/12/2020 2:30pm-5:50pm', 'Case attended site'),
    @('Lakes Entrance', 'Wyanga Winery  248 Baades Rd, Lakes Entrance VIC 3909', '29/12/2020 1:00pm-2:00pm', 'Case visited vineyard'),
    @('Leongatha', 'Coral Fish and Chips 53 Bair St, Leongatha VIC 3953', '26/12/2020 5:30pm-7:30pm', 'Case worked in store'),
    @('Leongatha', 'Coral Fish and Chips 53 Bair St, Leongatha VIC 3953', '27/12/2020 4:00pm-7:30pm', 'Case worked in store'),
    @('Melbourne', 'European Bier Cafe City  120 Exhibition Street Melbourne VIC 3000', '21/12/2020 8:00pm-9:30pm', 'Case attended cafe'),
    @('Melbourne', 'Fonda Mexican Flinders Lane  31 Flinders Lane Melbourne', '29/12/2020 6:00pm-7:00pm', 'Case attended restaurant'),
    @('Moorabbin', 'Minnie Miny Mo Cafe  8 Station Street, Moorabbin', '24/12/2020 11:00am-12:00pm', 'Case attended cafe'),
    @('Oakleigh', 'Melissa Oakleigh Restaurant  6 Eaton Mall, Oakleigh VIC 3166', '28/12/20 7:30pm-8:15pm', 'Case dined in restaurant'),
    @('Southbank', 'Left Bank Melbourne Restaurant and Cocktail Bar  1 Southbank Boulevard, Southbank', '25/12/2020 12:00pm-02:30pm', 'Case attended bar'),
    @('Southbank', 'Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank', '23/112/2020 1:00pm-1:30pm', 'Case attended restaurant')
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}

# The table is now one row shorter (old last data row was 18, new is 17),
# so clear out the now-stale trailing row.
$lastOldRow = 18
for ($r = $rowIndex; $r -le $lastOldRow; $r++) {
    $ws.Range("A$r`:D$r").Clear()
}
